# Add the new "news channel features" columns (F:K) to Sheet1, mirroring the
# existing A:E layout (header row styled "Bad", numeric data below for the
# rows that already have source data), then adjust the view to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:K1 (shared strings + "Bad" style, like B1:E1) ---
$headers = @("bbc", "cnn", "cnnibn", "ndtv", "timesnow", "features")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 6 + $i
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Style = "Bad"
}

# --- New numeric data for rows 2-6 (F:K) ---
$data = @(
    @(5319.97,            9657.02,            18029.900000000001, 4969.1400000000003, 25186,   520072.6),
    @(2611.6,             4904.2299999999996, 9459.1200000000008, 2815.15,             13141.6, 274190.53000000003),
    @(1731.61,            2567.88,            5286.9,             1587.57,             7120.98, 154592.4),
    @(1099.1300000000001, 1508.05,            3185.74,            977.59699999999998,  4229.33, 75923.199999999997),
    @(728.529,            1057.1099999999999, 2030.79,            645.07899999999995,  2733.85, 38885.5)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 2 + $r
    for ($c = 0; $c -lt 6; $c++) {
        $col = 6 + $c
        $ws.Cells.Item($row, $col).Value = $data[$r][$c]
    }
}

# --- Column J (10) gets an explicit width, like the other data columns ---
$ws.Columns.Item(10).ColumnWidth = 9.75

# --- View: zoom out to 145% and move the selection to K6 ---
$excel.ActiveWindow.Zoom = 145
$ws.Range("K6").Select() | Out-Null
